# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N - shifting the existing "Late" / "heading" / "Outstanding"
# columns one place to the right - and make this sheet the active tab
# with cell R6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet (this moves tabSelected /
# activeTab onto it and off whatever sheet was active before).
$ws.Activate()

# Insert a new blank column at N (pushes N:P -> O:Q).
$ws.Columns("N:N").Insert()

# The newly inserted column inherits the width of its left neighbour (M).
$ws.Columns("N:N").ColumnWidth = 10.14

# Update the selection to R6, as left by the author after the edit.
$ws.Range("R6").Select() | Out-Null
